$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new columns M, N, O ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Copy the header formatting (bold font + borders) from an existing header cell (L1)
# onto the three new header cells so they match the rest of the header row style.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2..41: populate the new columns ---
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"      # column M: renewd
    $ws.Cells.Item($r, 14).Value = 20120894      # column N: PlanID
    $ws.Cells.Item($r, 15).Value = 6             # column O: iteration
}
